$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Participant query text in B2 ---
# The original "simple" participant query is replaced by a more elaborate
# one using OPTIONAL MATCH / apoc.coll.sort, etc. Replacing the cell value
# drops the now-unused shared string and appends the new one at the end of
# the shared-strings table, which is exactly what the target workbook does
# (the Stat/Sample/File query strings keep their text but shift index).
$newParticipantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.platform in ['Illumina Next Seq 500']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Range("B2").Value = $newParticipantQuery

# --- Grow row 2 to fit the now-longer query text ---
$ws.Rows.Item(2).RowHeight = 279

# --- Update the view: scroll down a bit and move the active selection ---
$ws.Activate()
$ws.Range("B5").Select()
